# Updated symbol list on Tue Dec 20 14:34:59 UTC 2022 with GitHub Actions
#
# This script updates the "Price" column (D) for several coins with refreshed
# quotes, and inserts a new "One" entry at row 10 pushing the WazirX..CoinExToken
# block down by one row (row 18's "One" entry is replaced by the block shift).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Text
    )
    $rng = $ws.Range($Address)
    # Prefix with an apostrophe so Excel always stores the value as text,
    # even when the text looks like a number (e.g. "249.15").
    $rng.Value = "'" + $Text
    # Re-apply the default "Normal" style so the forced text entry does not
    # leave a stray quotePrefix-flavoured style behind on the cell.
    $rng.Style = "Normal"
}

# --- Simple price refreshes (column D only) ---
Set-TextValue "D2" "249.15"
Set-TextValue "D3" "22.91"
Set-TextValue "D4" "5.389"
Set-TextValue "D5" "0.05614"
Set-TextValue "D7" "6.358"
Set-TextValue "D8" "0.8158"
Set-TextValue "D9" "0.9174"

# --- Rows 10-18: a new "One" ranking entry appears at row 10, shifting the
#     existing WazirX..CoinExToken block down by one row. Columns B, C and E
#     move down to the next row (keeping the row's own rank prefix in E),
#     while column D receives a freshly refreshed price for each coin. ---

Set-TextValue "B10" "One"
Set-TextValue "C10" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D10" "0.01136"
Set-TextValue "E10" "9OneONE"

Set-TextValue "B11" "WazirX"
Set-TextValue "C11" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D11" "0.1423"
Set-TextValue "E11" "10WazirXWRX"

Set-TextValue "B12" "MandalaExchangeToken"
Set-TextValue "C12" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D12" "0.07484"
Set-TextValue "E12" "11MandalaExchangeTokenMDX"

Set-TextValue "B13" "LiechtensteinCryptoassetsExchange"
Set-TextValue "C13" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D13" "0.03193"
Set-TextValue "E13" "12LiechtensteinCryptoassetsExchangeLCX"

Set-TextValue "B14" "BitrueCoin"
Set-TextValue "C14" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D14" "0.03096"
Set-TextValue "E14" "13BitrueCoinBTR"

Set-TextValue "B15" "BitMartToken"
Set-TextValue "C15" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D15" "0.09323"
Set-TextValue "E15" "14BitMartTokenBMX"

Set-TextValue "B16" "MCDex"
Set-TextValue "C16" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D16" "3.570"
Set-TextValue "E16" "15MCDexMCB"

Set-TextValue "B17" "BitForexToken"
Set-TextValue "C17" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D17" "0.001590"
Set-TextValue "E17" "16BitForexTokenBF"

Set-TextValue "B18" "CoinExToken"
Set-TextValue "C18" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D18" "0.04719"
Set-TextValue "E18" "17CoinExTokenCET"

# --- Remaining scattered price refreshes (column D only) ---
Set-TextValue "D19" "0.006403"
Set-TextValue "D21" "0.001033"
Set-TextValue "D24" "2.164"
Set-TextValue "D25" "0.3250"
Set-TextValue "D28" "0.0003000"
Set-TextValue "D40" "0.03965"
Set-TextValue "D41" "0.006902"
Set-TextValue "D42" "0.1064"
Set-TextValue "D44" "0.007537"
Set-TextValue "D45" "0.00005573"
Set-TextValue "D48" "0.6753"
Set-TextValue "D49" "0.2199"
Set-TextValue "D50" "0.00002100"
